$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (rows 7-11) to append to the running log sheet
$data = @(
    @("20160401_015227", 422.647, "convert to lower, trim `"space`" and `",`", remove multiple spaces, convert unicode to ascii", "2 features: #ascii/(#ascii+#digit+#punctuation), #max_digit_skip_0_1", "Neuron Network", "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 300", 0.966, 0.996699669966997, "0 filters: ", 0.26530612244898),
    @("20160401_015929", 583.758, "convert to lower, trim `"space`" and `",`", remove multiple spaces, convert unicode to ascii", "2 features: #ascii/(#ascii+#digit+#punctuation), #max_digit_skip_0_1", "Neuron Network", "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 300", 0.960666666666667, 0.996699669966997, "0 filters: ", 0.255102040816327),
    @("20160401_020913", 533.151, "convert to lower, trim `"space`" and `",`", remove multiple spaces, convert unicode to ascii", "2 features: #ascii/(#ascii+#digit+#punctuation), #max_digit_skip_0_1", "Neuron Network", "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 300", 0.961333333333333, 0.996699669966997, "0 filters: ", 0.244897959183673),
    @("20160401_021806", 537.55, "convert to lower, trim `"space`" and `",`", remove multiple spaces, convert unicode to ascii", "2 features: #ascii/(#ascii+#digit+#punctuation), #max_digit_skip_0_1", "Neuron Network", "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 300", 0.956, 0.996699669966997, "0 filters: ", 0.23469387755102),
    @("20160401_022704", 792.787, "convert to lower, trim `"space`" and `",`", remove multiple spaces, convert unicode to ascii", "2 features: #ascii/(#ascii+#digit+#punctuation), #max_digit_skip_0_1", "Neuron Network", "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 300", 0.955333333333333, 0.996699669966997, "0 filters: ", 0.23469387755102)
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
    $ws.Cells.Item($row, 8).Value = $rowData[7]
    $ws.Cells.Item($row, 9).Value = $rowData[8]
    $ws.Cells.Item($row, 10).Value = $rowData[9]
}
